# Apply the registry update described by the commit:
# "Finalize v18 reinstatement, source audit sync, and full QA rebuild"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Sources": append new row 129 for the canonical legacy source S128
# ---------------------------------------------------------------------
$sources = $wb.Worksheets.Item("Sources")

$sources.Range("A129").Value = "S128"
$sources.Range("B129").Value = "Legacy v18 Reference Archive"
$sources.Range("C129").Value = "Internal Legacy Document"
$sources.Range("D129").Value = "sources/internal/S128_legacy_v18_reference.docx"
$sources.Range("F129").Value = "'2026-02-08"
$sources.Range("G129").Value = "'2026-02-08"
$sources.Range("H129").Value = "Codex"
$sources.Range("I129").Value = "Canonical legacy benchmark source copied from _output/archive/v18.docx into sources/internal to avoid circular source dependency on generated outputs. Used for reinstated legacy tables, investor map extensions, and Figure 44 context."

# ---------------------------------------------------------------------
# Sheet "Claims": sync source audit statuses and notes
# ---------------------------------------------------------------------
$claims = $wb.Worksheets.Item("Claims")

# Row 54 (C053) - deprecated from active section wording, mapped to S128
$claims.Range("D54").Value = "S128"
$claims.Range("E54").Value = "I.2 historical note (deprecated)"
$claims.Range("J54").Value = "[2026-02-08] Deprecated from active section wording; retained for audit trail and mapped to legacy v18 source S128."

# Row 55 (C054) - claim text rewritten, remains UNVERIFIED, location added
$claims.Range("C55").Value = "The specific threshold convention used in prior internal work (>5% R&D and >20% EBITDA as repeatable premium signal) remains unresolved and stays flagged [UNVERIFIED]."
$claims.Range("E55").Value = "I.2 paragraph 2"
$claims.Range("J55").Value = "[2026-02-08] Active unresolved claim in section text; threshold heuristic still lacks direct source binding."

# Row 57 (C056) - claim text rewritten, remains UNVERIFIED, location added
$claims.Range("C57").Value = "Specific monetization splits for selected overlays in the 11-segment matrix remain unresolved and are retained transparently [UNVERIFIED]."
$claims.Range("E57").Value = "I.3 paragraph 3"
$claims.Range("J57").Value = "[2026-02-08] Active unresolved claim in section text (overlay monetization splits)."

# Row 59 (C058) - claim text rewritten, remains UNVERIFIED, location added
$claims.Range("C59").Value = "Comparative concentration assumptions across validated versus commodity categories remain unresolved and are retained transparently [UNVERIFIED]."
$claims.Range("E59").Value = "I.3 paragraph 3"
$claims.Range("J59").Value = "[2026-02-08] Active unresolved claim in section text (concentration assumptions)."

# Row 71 (C070) - resolved from UNVERIFIED to S128
$claims.Range("C71").Value = "The specific legacy SOM point estimate is retained as a legacy benchmark from v18/v19 source material [S128]."
$claims.Range("D71").Value = "S128"
$claims.Range("E71").Value = "III.1 paragraph 3"
$claims.Range("J71").Value = "[2026-02-08] Resolved from UNVERIFIED to S128 after adding canonical legacy source artifact."

# Row 82 (C081) - resolved from UNVERIFIED to S128
$claims.Range("C82").Value = "Figure 44 Opportunity matrix is sourced to legacy v18 benchmark material [S128]."
$claims.Range("D82").Value = "S128"
$claims.Range("J82").Value = "[2026-02-08] Resolved from UNVERIFIED to S128."

# Row 89 (C088) - deprecated note, mapped to S128
$claims.Range("D89").Value = "S128"
$claims.Range("E89").Value = "III.2 legacy note (deprecated)"
$claims.Range("J89").Value = "[2026-02-08] Deprecated from active section text; retained for audit trail and mapped to S128."

# Row 90 (C089) - deprecated note, mapped to S128
$claims.Range("D90").Value = "S128"
$claims.Range("E90").Value = "III.2 legacy note (deprecated)"
$claims.Range("J90").Value = "[2026-02-08] Deprecated from active section text; retained for audit trail and mapped to S128."

# Row 91 (C090) - deprecated note, mapped to S128
$claims.Range("D91").Value = "S128"
$claims.Range("E91").Value = "III.2 legacy note (deprecated)"
$claims.Range("J91").Value = "[2026-02-08] Deprecated from active section text; retained for audit trail and mapped to S128."

# ---------------------------------------------------------------------
# Sheet "Figures": FIG-44 resolved from UNVERIFIED to S128, status Final
# ---------------------------------------------------------------------
$figures = $wb.Worksheets.Item("Figures")

$figures.Range("D45").Value = "S128"
$figures.Range("G45").Value = "Final"
$figures.Range("I45").Value = "[2026-02-08] Resolved from UNVERIFIED to S128 after adding canonical legacy source artifact in sources/internal."
